$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88
$ws.Range("A88").Value = 26
$ws.Range("B88").Value = 66.05
$ws.Range("C88").Value = 12
$ws.Range("D88").Value = 300
$ws.Range("E88").Value = 2000
$ws.Range("F88").Value = "HWBAQÓ ŚDJGIŃ ĄCĆEĘF KLŁMNO PRSTUV XYZŹŻ|"
$ws.Range("G88").Value = -1403.88365943578
$ws.Range("H88").Value = "HWBAQÓ ŚDJGIŃ ĄCĆEĘF KLŁMNO PRSTUV XYZŹŻ|"
$ws.Range("I88").Value = -1403.8837

# Row 89
$ws.Range("A89").Value = 46
$ws.Range("B89").Value = 155.07
$ws.Range("C89").Value = 12
$ws.Range("D89").Value = 300
$ws.Range("E89").Value = 2000
$ws.Range("F89").Value = "DUZXŚI TĆVSŹA ĄBCEĘF GHJKLŁ MNŃOÓP QRWYŻ|"
$ws.Range("G89").Value = -1403.88365943578
$ws.Range("H89").Value = "DUZXŚI TĆVSŹA ĄBCEĘF GHJKLŁ MNŃOÓP QRWYŻ|"
$ws.Range("I89").Value = -1403.8837

# Row 90
$ws.Range("A90").Value = 105
$ws.Range("B90").Value = 367.03
$ws.Range("C90").Value = 12
$ws.Range("D90").Value = 300
$ws.Range("E90").Value = 2000
$ws.Range("F90").Value = "TYIMCJ Ż|ŃĆPN AĄBDEĘ FGHKLŁ OÓQRSŚ UVWXZŹ"
$ws.Range("G90").Value = -1403.88365943578
$ws.Range("H90").Value = "TYIMCJ Ż|ŃĆPN AĄBDEĘ FGHKLŁ OÓQRSŚ UVWXZŹ"
$ws.Range("I90").Value = -1403.8837

# Row 91
$ws.Range("A91").Value = 15
$ws.Range("B91").Value = 53.42
$ws.Range("C91").Value = 12
$ws.Range("D91").Value = 300
$ws.Range("E91").Value = 2000
$ws.Range("F91").Value = "JUIRAF ŁNWHZK ĄBCĆDE ĘGLMŃO ÓPQSŚT VXYŹŻ|"
$ws.Range("G91").Value = -1403.88365943578
$ws.Range("H91").Value = "JUIRAF ŁNWHZK ĄBCĆDE ĘGLMŃO ÓPQSŚT VXYŹŻ|"
$ws.Range("I91").Value = -1403.8837

# Row 92
$ws.Range("A92").Value = 50
$ws.Range("B92").Value = 168.1
$ws.Range("C92").Value = 12
$ws.Range("D92").Value = 300
$ws.Range("E92").Value = 2000
$ws.Range("F92").Value = "OÓIEŚR ŹBTYLK AĄCĆDĘ FGHJŁM NŃPQSU VWXZŻ|"
$ws.Range("G92").Value = -1403.88365943578
$ws.Range("H92").Value = "OÓIEŚR ŹBTYLK AĄCĆDĘ FGHJŁM NŃPQSU VWXZŻ|"
$ws.Range("I92").Value = -1403.8837

# Row 93
$ws.Range("A93").Value = 9
$ws.Range("B93").Value = 34.72
$ws.Range("C93").Value = 12
$ws.Range("D93").Value = 300
$ws.Range("E93").Value = 2000
$ws.Range("F93").Value = "ŹVFCEĄ SYDŃAJ BĆĘGHI KLŁMNO ÓPQRŚT UWXZŻ|"
$ws.Range("G93").Value = -1403.88365943578
$ws.Range("H93").Value = "ŹVĘCEĄ SYDŃAJ BĆFGHI KLŁMNO ÓPQRŚT UWXZŻ|"
$ws.Range("I93").Value = -1471.1039

# Row 94
$ws.Range("A94").Value = 3
$ws.Range("B94").Value = 15.51
$ws.Range("C94").Value = 12
$ws.Range("D94").Value = 300
$ws.Range("E94").Value = 2000
$ws.Range("F94").Value = "VŻBTMW FŚXÓRQ AĄCĆDE ĘGHIJK LŁNŃOP SUYZŹ|"
$ws.Range("G94").Value = -1403.88365943578
$ws.Range("H94").Value = "VŻBTMW FŚXÓRQ AĄCĆDE ĘGHIJK LŁNŃOP SUYZŹ|"
$ws.Range("I94").Value = -1403.8837

# Row 95
$ws.Range("A95").Value = 18
$ws.Range("B95").Value = 54.78
$ws.Range("C95").Value = 12
$ws.Range("D95").Value = 300
$ws.Range("E95").Value = 2000
$ws.Range("F95").Value = "TŃĘ|ŻŹ FNMAEÓ ĄBCĆDG HIJKLŁ OPQRSŚ UVWXYZ"
$ws.Range("G95").Value = -1403.88365943578
$ws.Range("H95").Value = "T|ĘŃŻŹ FNMAEÓ ĄBCĆDG HIJKLŁ OPQRSŚ UVWXYZ"
$ws.Range("I95").Value = -1446.2486

# Row 96
$ws.Range("A96").Value = 19
$ws.Range("B96").Value = 55.58
$ws.Range("C96").Value = 12
$ws.Range("D96").Value = 300
$ws.Range("E96").Value = 2000
$ws.Range("F96").Value = "ŚJŻBĘL KŹINZP AĄCĆDE FGHŁMŃ OÓQRST UVWXY|"
$ws.Range("G96").Value = -1403.88365943578
$ws.Range("H96").Value = "OÓQRST UVWXY| ŚJŻBĘL KŹINZP AĄCĆDE FGHŁMŃ"
$ws.Range("I96").Value = -1403.8837

# Row 97
$ws.Range("A97").Value = 149
$ws.Range("B97").Value = 450.01
$ws.Range("C97").Value = 12
$ws.Range("D97").Value = 300
$ws.Range("E97").Value = 2000
$ws.Range("F97").Value = "RGAŃZĘ TPCMŻF ĄBĆDEH IJKLŁN OÓQSŚU VWXYŹ|"
$ws.Range("G97").Value = -1403.88365943578
$ws.Range("H97").Value = "OÓQSŚU VWXYŹ| RGAŃZĘ TPCMŻF ĄBHDEĆ IJKLŁN"
$ws.Range("I97").Value = -1461.2366

# Row 98
$ws.Range("A98").Value = 15
$ws.Range("B98").Value = 46.82
$ws.Range("C98").Value = 12
$ws.Range("D98").Value = 300
$ws.Range("E98").Value = 2000
$ws.Range("F98").Value = "ÓUBĄŹŃ MHTYAC ĆDEĘFG IJKLŁN OPQRSŚ VWXZŻ|"
$ws.Range("G98").Value = -1403.88365943578
$ws.Range("H98").Value = "ÓUŃĄŹB MHTYAC ĆDEĘFG IJKLŁN OPQRSŚ VWXZŻ|"
$ws.Range("I98").Value = -1469.9451

# Row 99
$ws.Range("A99").Value = 13
$ws.Range("B99").Value = 43.09
$ws.Range("C99").Value = 12
$ws.Range("D99").Value = 300
$ws.Range("E99").Value = 2000
$ws.Range("F99").Value = "VZXGĄB ĆCŚKSF ADEĘHI JLŁMNŃ OÓPQRT UWYŹŻ|"
$ws.Range("G99").Value = -1403.88365943578
$ws.Range("H99").Value = "VZXGĄB ĆCŚKSF ADEĘHI JLŁMNŃ OÓPQRT UWYŹŻ|"
$ws.Range("I99").Value = -1403.8837

# Row 100
$ws.Range("A100").Value = 21
$ws.Range("B100").Value = 95.46
$ws.Range("C100").Value = 12
$ws.Range("D100").Value = 300
$ws.Range("E100").Value = 4000
$ws.Range("F100").Value = "HŃJAŚI LV|OCĘ ĄBĆDEF GKŁMNÓ PQRSTU WXYZŹŻ"
$ws.Range("G100").Value = -1403.88365943578
$ws.Range("H100").Value = "HŃJAŚI LV|OCĘ ĄBĆDEF GKŁMNÓ PQRSTU WXYZŹŻ"
$ws.Range("I100").Value = -1403.8837

# Row 101
$ws.Range("A101").Value = 9
$ws.Range("B101").Value = 45.99
$ws.Range("C101").Value = 12
$ws.Range("D101").Value = 300
$ws.Range("E101").Value = 4000
$ws.Range("F101").Value = "XŻEŹĘL ŃŚKV|B AĄCĆDF GHIJŁM NOÓPQR STUWYZ"
$ws.Range("G101").Value = -1403.88365943578
$ws.Range("H101").Value = "XŻEŹĘL ŃŚKV|B AĄCĆDF GHIJŁM NOÓPQR STUWYZ"
$ws.Range("I101").Value = -1403.8837

# Row 102
$ws.Range("A102").Value = 26
$ws.Range("B102").Value = 114.93
$ws.Range("C102").Value = 12
$ws.Range("D102").Value = 300
$ws.Range("E102").Value = 4000
$ws.Range("F102").Value = "PFJKQO ZÓĄVŚW ABCĆDE ĘGHILŁ MNŃRST UXYŹŻ|"
$ws.Range("G102").Value = -1403.88365943578
$ws.Range("H102").Value = "PQJKFO ZÓĄVŚW ABCĆDE ĘGHILŁ MNŃRST UXYŹŻ|"
$ws.Range("I102").Value = -1450.9253
